{"js": "// Split the three long, single-run \"wall of text\" paragraphs (the\n// Portuguese \"Programa\", the italic English \"Programa\" and the\n// \"Bibliografia\" list) into multiple <w:t> runs separated by manual line\n// breaks (<w:br/>), one break between each numbered topic in the program\n// paragraphs and two breaks between each reference in the bibliography.\n//\n// Office.js has no direct \"insert <w:br/> at this offset\" primitive, but\n// Range.insertText(text, \"Replace\") maps \"\\v\" (vertical tab) onto a\n// w:br (manual line break) when it rewrites the backing run(s), so we\n// rebuild each paragraph's text with \"\\v\" at the desired break points and\n// replace the whole paragraph range in one shot. That keeps the run's\n// existing formatting (e.g. the italic rPr) untouched because the replace\n// happens on the paragraph's existing range/run rather than creating a\n// brand-new one.\n\nconst body = context.document.body;\n\nasync function splitParagraphByText(needle, segments, joiner) {\n  const results = body.search(needle, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not locate paragraph starting with: \" + needle);\n  }\n\n  const paragraph = results.items[0].paragraphs.getFirst();\n  const newText = segments.join(joiner);\n  paragraph.insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// 1) \"Programa\" (Portuguese) - one <w:br/> between each of the 11 topics.\nconst programaPt = [\n  \"1 - Introdu\u00e7\u00e3o ao Laborat\u00f3rio: No\u00e7\u00f5es Elementares de Seguran\u00e7a; Equipamentos B\u00e1sicos de Laborat\u00f3rio; Equipamentos de Prote\u00e7\u00e3o Individual.\",\n  \"2 - Pesos e medidas (Tratamento de dados experimentais): Cuidados Gerais com Balan\u00e7as; T\u00e9cnicas de Determina\u00e7\u00e3o de massa; Exatid\u00e3o e precis\u00e3o; Unidades; Algarismos Significativos; Propaga\u00e7\u00e3o de Erros.\",\n  \"3 - T\u00e9cnicas de Separa\u00e7\u00e3o de Misturas: Filtra\u00e7\u00e3o simples; Filtra\u00e7\u00e3o a v\u00e1cuo e Decanta\u00e7\u00e3o.\",\n  \"4 - Fen\u00f4menos f\u00edsicos: Constru\u00e7\u00e3o do Diagrama da mudan\u00e7a do estado f\u00edsico da \u00e1gua.\",\n  \"5 - Miscibilidade e solubilidade: Influ\u00eancia das for\u00e7as intermoleculares na miscibilidade de l\u00edquidos.\",\n  \"6 - Rea\u00e7\u00f5es qu\u00edmicas: Aspectos qualitativos.\",\n  \"7 - Solu\u00e7\u00f5es: Preparo e padroniza\u00e7\u00e3o de solu\u00e7\u00f5es.\",\n  \"8 - Titrimetria: Realiza\u00e7\u00e3o de Titula\u00e7\u00f5es \u00c1cido-Base; Retrotitula\u00e7\u00e3o.\",\n  \"9 - Equil\u00edbrio Qu\u00edmico - Princ\u00edpio de Le Chatelier\",\n  \"10 - Fatores que alteram a velocidade das rea\u00e7\u00f5es: Velocidade das rea\u00e7\u00f5es\",\n  \"11 \u2013 Termoqu\u00edmica: Entalpia de decomposi\u00e7\u00e3o do H2O2\",\n];\nawait splitParagraphByText(\n  \"1 - Introdu\u00e7\u00e3o ao Laborat\u00f3rio: No\u00e7\u00f5es Elementares de Seguran\u00e7a\",\n  programaPt,\n  \"\\v\"\n);\n\n// 2) \"Programa\" (italic English) - only the final topic (11) is split out\n// with a single <w:br/>; everything before stays as one run of text.\nconst programaEn = [\n  \"1 - Introduction to the Chemistry Laboratory: Elementary notion of security, Laboratory basic equipment; Individual protection equipment. 2 - Weights and measures (experimental data treatment): General care with scales, Determination of mass techniques. Accuracy and precision, units, significant digits and error propagation. 3 - Methods for separating mixtures: Simple filtration; Vacuum filtration and Decantation. 4 - Physical phenomena: Water state changes. 5 - Miscibility and solubility: Intermolecular forces influence on the liquids miscibility. 6 - Chemical reactions: Qualitative aspects. 7 - Solutions: Preparation and standardization of solutions. 8 - Titrimetry: Acid-Base Titrations and return-titration. 9 - Chemical Equilibrium - Le Chatelier's Principle. 10 - Factors that change the speed of reactions: Speed of reactions.\",\n  \"11 \u2013 Thermochemistry: Enthalpy of decomposition of H2O2\",\n];\nawait splitParagraphByText(\n  \"1 - Introduction to the Chemistry Laboratory: Elementary notion\",\n  programaEn,\n  \"\\v\"\n);\n\n// 3) \"Bibliografia\" - two <w:br/> (a blank line) between each reference.\nconst bibliografia = [\n  \"ASSUMP\u00c7\u00c3O, R. M. V.; MORITA, T. Manual de solu\u00e7\u00f5es reagentes e solventes: padroniza\u00e7\u00e3o, prepara\u00e7\u00e3o, purifica\u00e7\u00e3o. S\u00e3o Paulo: Editora Edgard Blucher, 2\u00aa Ed, 2007.\",\n  \"ATKINS, P.; JONES, L. Princ\u00edpios de Qu\u00edmica, Questionando a vida e o meio ambiente, Bookman, Porto Alegre, 5\u00aa Ed, 2011.\",\n  \"BACCAN, N.; ANDRADE, J. C. O.; GODINHO, E. S.; BARONE, J. S. Qu\u00edmica anal\u00edtica quantitativa elementar. 3.ed. S\u00e3o Paulo: Edgard Blucher, 2001.\",\n  \"BRADY,  J.E.; RUSSELL, J. W.; HOLUM, J.R. Qu\u00edmica - a Mat\u00e9ria e Suas Transforma\u00e7\u00f5es, 5\u00aa ed, Volume 1 e 2, LTC Editora, Rio de Janeiro, 2012.\",\n  \"BROWN, T. E et al. Qu\u00edmica a Ci\u00eancia Central. 9 ed. S\u00e3o Paulo. Pearson Prentice Hall, 2005-2007.\",\n  \"CONSTANTINO, M.G; SILVA, G. V. J. da; DONATE P. M. Fundamentos de qu\u00edmica experimental, S\u00e3o Paulo: EDUSP, 2004.\",\n  \"KOTZ, J.; TREICHEL, P.; WEAVER, G. Qu\u00edmica Geral e Rea\u00e7\u00f5es Qu\u00edmicas, Vol. 1 e 2, Cengage Learning, S\u00e3o Paulo, 2023.\",\n  \"MAHAN, B. M.; MYERS, R. J. Qu\u00edmica um curso universit\u00e1rio. S\u00e3o Paulo: Ed. Edgard Blucher Ltda, 1993.\",\n  \"SILVA, R. R.; BOCCHI, N.; ROCHA FILHO, R. P. Introdu\u00e7\u00e3o a qu\u00edmica experimental. S\u00e3o Paulo: a: EDUFSCAR, 2019\",\n];\nawait splitParagraphByText(\n  \"ASSUMP\u00c7\u00c3O, R. M. V.; MORITA, T. Manual de solu\u00e7\u00f5es\",\n  bibliografia,\n  \"\\v\\v\"\n);\n", "ps1": "# Split the three long, single-run \"wall of text\" paragraphs (the\n# Portuguese \"Programa\", the italic English \"Programa\" and the\n# \"Bibliografia\" list) into multiple runs of text separated by manual\n# line breaks (Word's vertical-tab / Chr(11) character, which saves as\n# <w:br/> in the OOXML), one break between each numbered topic in the\n# program paragraphs and two breaks (a blank line) between each\n# reference in the bibliography.\n#\n# Assigning Range.Text with Chr(11) embedded is exactly what Word does\n# when a user presses Shift+Enter, so it rewrites the backing run(s) as\n# <w:t>...</w:t><w:br/><w:t>...</w:t> while keeping the run's existing\n# formatting (e.g. the italic rPr) intact, since it edits the paragraph's\n# existing range/run instead of creating a brand-new one.\n\n$d = $word.ActiveDocument\n$LF = [char]11\n\nfunction Get-ParagraphRangeByText($doc, $needle) {\n    $rng = $doc.Content\n    $found = $rng.Find.Execute($needle)\n    if (-not $found) {\n        throw \"Could not locate paragraph starting with: $needle\"\n    }\n    $rng.Expand(4) | Out-Null   # wdParagraph\n    return $rng\n}\n\n# 1) \"Programa\" (Portuguese) - one manual line break between each of the\n#    11 topics.\n$programaPt = @(\n    \"1 - Introdu\u00e7\u00e3o ao Laborat\u00f3rio: No\u00e7\u00f5es Elementares de Seguran\u00e7a; Equipamentos B\u00e1sicos de Laborat\u00f3rio; Equipamentos de Prote\u00e7\u00e3o Individual.\",\n    \"2 - Pesos e medidas (Tratamento de dados experimentais): Cuidados Gerais com Balan\u00e7as; T\u00e9cnicas de Determina\u00e7\u00e3o de massa; Exatid\u00e3o e precis\u00e3o; Unidades; Algarismos Significativos; Propaga\u00e7\u00e3o de Erros.\",\n    \"3 - T\u00e9cnicas de Separa\u00e7\u00e3o de Misturas: Filtra\u00e7\u00e3o simples; Filtra\u00e7\u00e3o a v\u00e1cuo e Decanta\u00e7\u00e3o.\",\n    \"4 - Fen\u00f4menos f\u00edsicos: Constru\u00e7\u00e3o do Diagrama da mudan\u00e7a do estado f\u00edsico da \u00e1gua.\",\n    \"5 - Miscibilidade e solubilidade: Influ\u00eancia das for\u00e7as intermoleculares na miscibilidade de l\u00edquidos.\",\n    \"6 - Rea\u00e7\u00f5es qu\u00edmicas: Aspectos qualitativos.\",\n    \"7 - Solu\u00e7\u00f5es: Preparo e padroniza\u00e7\u00e3o de solu\u00e7\u00f5es.\",\n    \"8 - Titrimetria: Realiza\u00e7\u00e3o de Titula\u00e7\u00f5es \u00c1cido-Base; Retrotitula\u00e7\u00e3o.\",\n    \"9 - Equil\u00edbrio Qu\u00edmico - Princ\u00edpio de Le Chatelier\",\n    \"10 - Fatores que alteram a velocidade das rea\u00e7\u00f5es: Velocidade das rea\u00e7\u00f5es\",\n    \"11 \u2013 Termoqu\u00edmica: Entalpia de decomposi\u00e7\u00e3o do H2O2\"\n)\n$rngPt = Get-ParagraphRangeByText $d \"Introdu\u00e7\u00e3o ao Laborat\u00f3rio: No\u00e7\u00f5es Elementares\"\n$rngPt.Text = ($programaPt -join $LF)\n\n# 2) \"Programa\" (italic English) - only the final topic (11) is split out\n#    with a single manual line break; everything before stays as one run\n#    of text.\n$programaEn = @(\n    \"1 - Introduction to the Chemistry Laboratory: Elementary notion of security, Laboratory basic equipment; Individual protection equipment. 2 - Weights and measures (experimental data treatment): General care with scales, Determination of mass techniques. Accuracy and precision, units, significant digits and error propagation. 3 - Methods for separating mixtures: Simple filtration; Vacuum filtration and Decantation. 4 - Physical phenomena: Water state changes. 5 - Miscibility and solubility: Intermolecular forces influence on the liquids miscibility. 6 - Chemical reactions: Qualitative aspects. 7 - Solutions: Preparation and standardization of solutions. 8 - Titrimetry: Acid-Base Titrations and return-titration. 9 - Chemical Equilibrium - Le Chatelier's Principle. 10 - Factors that change the speed of reactions: Speed of reactions.\",\n    \"11 \u2013 Thermochemistry: Enthalpy of decomposition of H2O2\"\n)\n$rngEn = Get-ParagraphRangeByText $d \"Introduction to the Chemistry Laboratory: Elementary notion\"\n$rngEn.Text = ($programaEn -join $LF)\n\n# 3) \"Bibliografia\" - two manual line breaks (a blank line) between each\n#    reference.\n$bibliografia = @(\n    \"ASSUMP\u00c7\u00c3O, R. M. V.; MORITA, T. Manual de solu\u00e7\u00f5es reagentes e solventes: padroniza\u00e7\u00e3o, prepara\u00e7\u00e3o, purifica\u00e7\u00e3o. S\u00e3o Paulo: Editora Edgard Blucher, 2\u00aa Ed, 2007.\",\n    \"ATKINS, P.; JONES, L. Princ\u00edpios de Qu\u00edmica, Questionando a vida e o meio ambiente, Bookman, Porto Alegre, 5\u00aa Ed, 2011.\",\n    \"BACCAN, N.; ANDRADE, J. C. O.; GODINHO, E. S.; BARONE, J. S. Qu\u00edmica anal\u00edtica quantitativa elementar. 3.ed. S\u00e3o Paulo: Edgard Blucher, 2001.\",\n    \"BRADY,  J.E.; RUSSELL, J. W.; HOLUM, J.R. Qu\u00edmica - a Mat\u00e9ria e Suas Transforma\u00e7\u00f5es, 5\u00aa ed, Volume 1 e 2, LTC Editora, Rio de Janeiro, 2012.\",\n    \"BROWN, T. E et al. Qu\u00edmica a Ci\u00eancia Central. 9 ed. S\u00e3o Paulo. Pearson Prentice Hall, 2005-2007.\",\n    \"CONSTANTINO, M.G; SILVA, G. V. J. da; DONATE P. M. Fundamentos de qu\u00edmica experimental, S\u00e3o Paulo: EDUSP, 2004.\",\n    \"KOTZ, J.; TREICHEL, P.; WEAVER, G. Qu\u00edmica Geral e Rea\u00e7\u00f5es Qu\u00edmicas, Vol. 1 e 2, Cengage Learning, S\u00e3o Paulo, 2023.\",\n    \"MAHAN, B. M.; MYERS, R. J. Qu\u00edmica um curso universit\u00e1rio. S\u00e3o Paulo: Ed. Edgard Blucher Ltda, 1993.\",\n    \"SILVA, R. R.; BOCCHI, N.; ROCHA FILHO, R. P. Introdu\u00e7\u00e3o a qu\u00edmica experimental. S\u00e3o Paulo: a: EDUFSCAR, 2019\"\n)\n$rngBib = Get-ParagraphRangeByText $d \"ASSUMP\u00c7\u00c3O, R. M. V.; MORITA, T. Manual de solu\u00e7\u00f5es\"\n$doubleBreak = \"$LF$LF\"\n$rngBib.Text = ($bibliografia -join $doubleBreak)\n"}
